$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly indexing bug-fix: recompute B:G for rows 2-11 (quarters Q0-Q9)
# with corrected values (N count shifted by one quarter as well).

$values = @{
    2  = @(0.5209010662194005, 0.5284173549292553, 0.3633847247877238, 0.6028140051356835, 0.3140407230176452, 15)
    3  = @(0.3491358409449706, 0.3876299128740889, 0.2064515430844597, 0.4543693905672561, 0.3017614833873206, 14)
    4  = @(0.2766139601088831, 0.3205523163110667, 0.1530262537301577, 0.3911857023590685, 0.287900825001641,  13)
    5  = @(0.4245700094424397, 0.4393557528778729, 0.2414648368369625, 0.4913907170846459, 0.2583974611168742, 12)
    6  = @(0.3793378286036473, 0.4013254374874798, 0.2006046709917273, 0.4478891280124215, 0.2497563433833,    11)
    7  = @(0.3288324600564539, 0.3601674043196504, 0.1589310008397884, 0.3986615116107754, 0.2375808962882351, 10)
    8  = @(0.3582448599611742, 0.3844285013741874, 0.1824260236518181, 0.427113595723454,  0.2466728085108265, 9)
    9  = @(0.3532650510030724, 0.381336062820872,  0.1830290214229363, 0.4278189119509986, 0.257976355866154,  8)
    10 = @(0.3173063146244254, 0.3418763137882747, 0.1463084296434522, 0.3825028491965154, 0.2307148046110977, 7)
    11 = @(0.3606037648954714, 0.379698193710172,  0.1799938953220499, 0.4242568742189688, 0.2448480836729456, 6)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowValues[0]
    $ws.Cells.Item($row, 3).Value = $rowValues[1]
    $ws.Cells.Item($row, 4).Value = $rowValues[2]
    $ws.Cells.Item($row, 5).Value = $rowValues[3]
    $ws.Cells.Item($row, 6).Value = $rowValues[4]
    $ws.Cells.Item($row, 7).Value = $rowValues[5]
}
